$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 2015
$ws.Range("C2").Value = "S-1"
$ws.Range("D2").Value = 100

# Row 3
$ws.Range("A3").Value = 2015
$ws.Range("C3").Value = "S-2"
$ws.Range("D3").Value = 120

# Row 4
$ws.Range("A4").Value = 2016
$ws.Range("C4").Value = "S-1"
$ws.Range("D4").Value = 200

# Row 5
$ws.Range("A5").Value = 2016
$ws.Range("C5").Value = "S-2"
$ws.Range("D5").Value = 150

# Row 6
$ws.Range("A6").Value = 2017
$ws.Range("C6").Value = "S-1"
$ws.Range("D6").Value = 280

# Row 7
$ws.Range("A7").Value = 2017
$ws.Range("C7").Value = "S-2"
$ws.Range("D7").Value = 250
